# Natmi following Dr Hou advice
# Update the Gdnf-Gfra2 LR-pairs sheet: recompute edge weights (E..T) for the
# existing Sending cluster -> Target cluster rows, and insert two new target
# cluster rows (M1, M2) before the self-referencing "sCs" row, which moves
# from row 4 to row 6.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# Row 2: sCs -> ECs (values updated in place)
# ---------------------------------------------------------------------
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 1.995314333333334
$ws.Range("H2").Value = 5.985943000000001
$ws.Range("I2").Value = 1
$ws.Range("J2").Value = 1
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 0.3696466666666667
$ws.Range("N2").Value = 1.10894
$ws.Range("O2").Value = 0.06847815280531702
$ws.Range("P2").Value = 0.06847815280531701
$ws.Range("Q2").Value = 0.737561292268889
$ws.Range("R2").Value = 6.638051630420001
$ws.Range("S2").Value = 0.06847815280531702
$ws.Range("T2").Value = 0.06847815280531701

# ---------------------------------------------------------------------
# Row 3: sCs -> FAPs (values updated in place)
# ---------------------------------------------------------------------
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 1.995314333333334
$ws.Range("H3").Value = 5.985943000000001
$ws.Range("I3").Value = 1
$ws.Range("J3").Value = 1
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 3.30984
$ws.Range("N3").Value = 9.92952
$ws.Range("O3").Value = 0.6131577793599756
$ws.Range("P3").Value = 0.6131577793599756
$ws.Range("Q3").Value = 6.604171193040001
$ws.Range("R3").Value = 59.43754073736001
$ws.Range("S3").Value = 0.6131577793599756
$ws.Range("T3").Value = 0.6131577793599756

# ---------------------------------------------------------------------
# Insert two new rows before the current row 4 (sCs -> sCs), shifting it
# down to row 6, then fill rows 4 (M1) and 5 (M2) and row 6 (sCs, values
# updated).
# ---------------------------------------------------------------------
$ws.Range("A4:A5").EntireRow.Insert()

# Row 4: sCs -> M1
$ws.Range("A4").Value = "sCs"
$ws.Range("B4").Value = "Gdnf"
$ws.Range("C4").Value = "Gfra2"
$ws.Range("D4").Value = "M1"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 1.995314333333334
$ws.Range("H4").Value = 5.985943000000001
$ws.Range("I4").Value = 1
$ws.Range("J4").Value = 1
$ws.Range("K4").Value = 1
$ws.Range("L4").Value = 0.3333333333333333
$ws.Range("M4").Value = 0.027417
$ws.Range("N4").Value = 0.082251
$ws.Range("O4").Value = 0.005079081416839622
$ws.Range("P4").Value = 0.005079081416839621
$ws.Range("Q4").Value = 0.05470553307700001
$ws.Range("R4").Value = 0.4923497976930001
$ws.Range("S4").Value = 0.005079081416839622
$ws.Range("T4").Value = 0.005079081416839621

# Row 5: sCs -> M2
$ws.Range("A5").Value = "sCs"
$ws.Range("B5").Value = "Gdnf"
$ws.Range("C5").Value = "Gfra2"
$ws.Range("D5").Value = "M2"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 1.995314333333334
$ws.Range("H5").Value = 5.985943000000001
$ws.Range("I5").Value = 1
$ws.Range("J5").Value = 1
$ws.Range("K5").Value = 1
$ws.Range("L5").Value = 0.3333333333333333
$ws.Range("M5").Value = 0.1051533333333333
$ws.Range("N5").Value = 0.31546
$ws.Range("O5").Value = 0.0194799701372169
$ws.Range("P5").Value = 0.0194799701372169
$ws.Range("Q5").Value = 0.2098139531977778
$ws.Range("R5").Value = 1.88832557878
$ws.Range("S5").Value = 0.0194799701372169
$ws.Range("T5").Value = 0.0194799701372169

# Row 6: sCs -> sCs (previously row 4, values updated)
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 1.995314333333334
$ws.Range("H6").Value = 5.985943000000001
$ws.Range("I6").Value = 1
$ws.Range("J6").Value = 1
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 1.585966333333333
$ws.Range("N6").Value = 4.757899
$ws.Range("O6").Value = 0.2938050162806509
$ws.Range("P6").Value = 0.2938050162806509
$ws.Range("Q6").Value = 3.164501357084112
$ws.Range("R6").Value = 28.480512213757
$ws.Range("S6").Value = 0.2938050162806509
$ws.Range("T6").Value = 0.2938050162806509
